# Matcher.py-style edit: add a "Matcher" lookup column (value 10) immediately
# to the left of each sheet's numeric "C" data column, pushing the old data
# one column to the right (B->C inherits the left-aligned style, the shifted
# data keeps the original style).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "A": every data row (1-11) gets the new column.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("A")
$wsA.Columns.Item(3).Insert(-4161)   # xlShiftToRight: old col C -> col D
for ($r = 1; $r -le 11; $r++) {
    $wsA.Cells.Item($r, 3).Value = 10
}

# Make "A" the active sheet / tab, with C3 selected.
$wsA.Activate()
$wsA.Range("C3").Select()

# ---------------------------------------------------------------------------
# Sheet "A Sorted by Matcher": only the rows that came from sheet "A" (the
# ones whose column A says "AAAA") get the new column; the "BBBBB" rows are
# left completely untouched.
# ---------------------------------------------------------------------------
$wsASorted = $wb.Worksheets.Item("A Sorted by Matcher")
$aSortedRows = @(1, 3, 5, 7, 10, 13, 16, 18, 20, 22, 24)
foreach ($r in $aSortedRows) {
    $wsASorted.Cells.Item($r, 4).Value = $wsASorted.Cells.Item($r, 3).Value2
    $wsASorted.Cells.Item($r, 3).Value = 10
}

# ---------------------------------------------------------------------------
# Sheet "B Sorted by Matcher": only the rows that came from sheet "A" (rows
# 20, 23, 26) get the new column; all the "BBBBB" rows stay untouched.
# ---------------------------------------------------------------------------
$wsBSorted = $wb.Worksheets.Item("B Sorted by Matcher")
$bSortedRows = @(20, 23, 26)
foreach ($r in $bSortedRows) {
    $wsBSorted.Cells.Item($r, 4).Value = $wsBSorted.Cells.Item($r, 3).Value2
    $wsBSorted.Cells.Item($r, 3).Value = 10
}

Write-Output "Matcher column inserted."
